$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2022" column (S) to the table, mirroring the existing "2021"
# column (R) formatting (header style in row 4, data-row style below), then
# overwrite with the new year header / 0-value data.
$null = $ws.Range("R4:R14").Copy($ws.Range("S4:S14"))

$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("S9").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("S11").Value = 0
$ws.Range("S12").Value = 0
$ws.Range("S13").Value = 0
$ws.Range("S14").Value = 0

# Match the author's final selection recorded in the workbook view.
$null = $ws.Range("R17").Select()
